$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '41.524.67'
Set-TextValue 'D3' '2.467.51'
Set-TextValue 'E3' '  -0.67%  '
Set-TextValue 'D4' '0.998'
Set-TextValue 'E4' '  -0.34%  '
Set-TextValue 'D5' '314.96'
Set-TextValue 'E5' '  +0.21%  '
Set-TextValue 'D6' '91.71'
Set-TextValue 'E6' '  -2.86%  '
Set-TextValue 'E7' '  -0.10%  '
Set-TextValue 'E8' '  -0.21%  '
Set-TextValue 'E9' '  +3.27%  '
Set-TextValue 'D10' '32.41'
Set-TextValue 'E10' '  -3.35%  '
Set-TextValue 'D11' '0.0794'
Set-TextValue 'E11' '  +1.68%  '
Set-TextValue 'E12' '  +0.81%  '
Set-TextValue 'D13' '2.849.71'
Set-TextValue 'E13' '  -0.54%  '
Set-TextValue 'D14' '6.85'
Set-TextValue 'E14' '  -1.07%  '
Set-TextValue 'D15' '16.03'
Set-TextValue 'E15' '  +3.52%  '
Set-TextValue 'D16' '2.478.99'
Set-TextValue 'E16' '  +1.94%  '
Set-TextValue 'E17' '  -2.13%  '
Set-TextValue 'D18' '41.541.18'
Set-TextValue 'E18' '  +0.20%  '
Set-TextValue 'D19' '6.50'
Set-TextValue 'E19' '  +2.61%  '
Set-TextValue 'D20' '0.0₃0944'
Set-TextValue 'E20' '  +2.08%  '
Set-TextValue 'D21' '71.17'
Set-TextValue 'E21' '  +3.15%  '
Set-TextValue 'D22' '11.08'
Set-TextValue 'E22' '  -1.46%  '
Set-TextValue 'D23' '236.65'
Set-TextValue 'E23' '  -0.21%  '
Set-TextValue 'D24' '2.73'
Set-TextValue 'E24' '  -1.45%  '
Set-TextValue 'B25' 'ImmutableX'
Set-TextValue 'C25' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D25' '1.91'
Set-TextValue 'E25' '  +0.17%  '
Set-TextValue 'B26' 'Dai'
Set-TextValue 'C26' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D26' '1.00'
Set-TextValue 'E26' '  -0.12%  '
Set-TextValue 'D27' '24.77'
Set-TextValue 'E27' '  +2.79%  '
Set-TextValue 'E28' '  -0.99%  '
Set-TextValue 'D29' '9.68'
Set-TextValue 'E29' '  -1.37%  '
Set-TextValue 'D30' '35.31'
Set-TextValue 'E30' '  -4.16%  '
Set-TextValue 'D31' '155.65'
Set-TextValue 'E31' '  +2.03%  '
Set-TextValue 'D32' '5.44'
Set-TextValue 'E32' '  -1.12%  '
Set-TextValue 'E33' '  +0.16%  '
Set-TextValue 'D34' '0.0760'
Set-TextValue 'E34' '  +0.29%  '
Set-TextValue 'D35' '17.17'
Set-TextValue 'E35' '  -5.73%  '
Set-TextValue 'B36' 'ApeXProtocol'
Set-TextValue 'C36' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D36' '2.37'
Set-TextValue 'E36' '  -3.62%  '
Set-TextValue 'B37' 'LidoDAOToken'
Set-TextValue 'C37' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D37' '2.88'
Set-TextValue 'E37' '  -6.75%  '
Set-TextValue 'E38' '  -0.24%  '
Set-TextValue 'B39' 'Kaspa'
Set-TextValue 'C39' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D39' '0.103'
Set-TextValue 'E39' '  +1.34%  '
Set-TextValue 'B40' 'ARBITRUM'
Set-TextValue 'C40' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D40' '1.78'
Set-TextValue 'E40' '  -5.27%  '
Set-TextValue 'D41' '4.00'
Set-TextValue 'E41' '  -3.86%  '
Set-TextValue 'E42' '  -0.34%  '
Set-TextValue 'D43' '1.941.72'
Set-TextValue 'E43' '  -2.45%  '
Set-TextValue 'E44' '  -1.18%  '
Set-TextValue 'D45' '18.77'
Set-TextValue 'E45' '  -5.67%  '
Set-TextValue 'D46' '2.91'
Set-TextValue 'E46' '  -3.70%  '
Set-TextValue 'D47' '9.07'
Set-TextValue 'E47' '  +2.14%  '
Set-TextValue 'D48' '2.707.91'
Set-TextValue 'E48' '  -0.78%  '
Set-TextValue 'D49' '97.15'
Set-TextValue 'E49' '  -0.20%  '
Set-TextValue 'D50' '67.16'
Set-TextValue 'E50' '  -3.55%  '
Set-TextValue 'B51' 'MultiversX'
Set-TextValue 'C51' 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue 'D51' '52.51'
Set-TextValue 'E51' '  +2.94%  '
